$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.536887325674664
$ws.Range("C2").Value = 0.178962441891555
$ws.Range("D2").Value = 8.67387631321838
$ws.Range("F2").Value = 0.0001

$ws.Range("B3").Value = 2.97105823299686
$ws.Range("C3").Value = 0.0206323488402559
